# Generate Report for Handoff
#
# The localization report is being regenerated: the "Status" / Overview
# handoff-state text flips from the previous handback state to
# "Ready for handoff", and the associated timestamps are refreshed to the
# new handoff run. The three report tables (Overview, zh-cn, de-de) also
# get their "Status" columns narrowed to their new auto-sized width.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status text: "Handed back: in sync with en-US" -> "Ready for handoff"
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value     = "Ready for handoff"
$dede.Range("C2").Value     = "Ready for handoff"

# --- Refreshed handoff timestamps for this run
$overview.Range("G2").Value = "2016-08-30 01:01:20"
$zhcn.Range("H2").Value     = "2016-08-30 01:01:16"
$dede.Range("H2").Value     = "2016-08-30 01:01:20"

# --- Narrow the Status columns to their new (re-sized) width.
# ColumnWidth is quantized by Excel to whole pixels, so this value is the
# closest representable width to the regenerated layout's target.
$overview.Columns.Item(5).ColumnWidth = 16.33
$overview.Columns.Item(6).ColumnWidth = 16.33
$zhcn.Columns.Item(3).ColumnWidth     = 16.33
$dede.Columns.Item(3).ColumnWidth     = 16.33
